$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 198; this shifts existing rows 198-211 down to 199-212
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new weekly price record
$ws.Cells.Item(198, 1).Value = 9
$ws.Cells.Item(198, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(198, 3).Value = "Metropolitana"
$ws.Cells.Item(198, 4).Value = 44585
$ws.Cells.Item(198, 5).Value = 13
$ws.Cells.Item(198, 6).Value = 100112030
$ws.Cells.Item(198, 7).Value = "Poroto granado"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 52
$ws.Cells.Item(198, 11).Value = 24000
$ws.Cells.Item(198, 12).Value = 25000
$ws.Cells.Item(198, 13).Value = 24500
$ws.Cells.Item(198, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(198, 15).Value = "Región del Maule"
$ws.Cells.Item(198, 16).Value = 980
$ws.Cells.Item(198, 17).Value = 25
$ws.Cells.Item(198, 18).Value = "Hortaliza"
